$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from the last existing header cell (AB1) into the new header cells
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill season record (Wins/Losses/Ties) for every data row
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 29).Value = 47
    $ws.Cells.Item($row, 30).Value = 70
    $ws.Cells.Item($row, 31).Value = 0
}
